$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record is inserted as row 261; every existing record from the
# old row 261 onward shifts down by one row (the former last row, 283,
# becomes the new row 284).
$ws.Rows.Item(261).Insert()

# Populate the newly inserted row 261 with this week's data (same
# Mercado/Region/Categoria/etc. as the rest of the block).
$ws.Cells.Item(261, 1).Value = 5
$ws.Cells.Item(261, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(261, 3).Value = "Maule"
$ws.Cells.Item(261, 4).Value = 44769
$ws.Cells.Item(261, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(261, 5).Value = 7
$ws.Cells.Item(261, 6).Value = 100112009
$ws.Cells.Item(261, 7).Value = "Acelga"
$ws.Cells.Item(261, 8).Value = "Sin especificar"
$ws.Cells.Item(261, 9).Value = "Primera"
$ws.Cells.Item(261, 10).Value = 500
$ws.Cells.Item(261, 11).Value = 4000
$ws.Cells.Item(261, 12).Value = 4000
$ws.Cells.Item(261, 13).Value = 4000
$ws.Cells.Item(261, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(261, 15).Value = "Región del Maule"
$ws.Cells.Item(261, 16).Value = 1000
$ws.Cells.Item(261, 17).Value = 4
$ws.Cells.Item(261, 18).Value = "Hortaliza"
